$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (same layout/
#    styles as the other quarterly fund-holding sheets) and placing it right
#    after "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($null, $totalSheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Update the existing holding (fund 010010) with the 2022-Q3 figures.
$q3Sheet.Cells.Item(2, 4).Value = "'6.96"
$q3Sheet.Cells.Item(2, 5).Value = "'90.21"
$q3Sheet.Cells.Item(2, 6).Value = "'3.83"
$q3Sheet.Cells.Item(2, 7).Value = "'0.2666"
$q3Sheet.Cells.Item(2, 8).Value = 8

# Clone the formatting of row 2 down into rows 3 and 4 for the two new
# holdings, then fill in their data.
$q3Sheet.Range("A2:H2").Copy()
$q3Sheet.Range("A3:H3").PasteSpecial(-4122)
$q3Sheet.Range("A4:H4").PasteSpecial(-4122)

$q3Sheet.Cells.Item(3, 1).Value = 1
$q3Sheet.Cells.Item(3, 2).Value = "'004403"
$q3Sheet.Cells.Item(3, 3).Value = "平安股息精选沪港深股票A"
$q3Sheet.Cells.Item(3, 4).Value = "'0.08"
$q3Sheet.Cells.Item(3, 5).Value = "'91.74"
$q3Sheet.Cells.Item(3, 6).Value = "'2.41"
$q3Sheet.Cells.Item(3, 7).Value = "'0.0019"
$q3Sheet.Cells.Item(3, 8).Value = 10

$q3Sheet.Cells.Item(4, 1).Value = 2
$q3Sheet.Cells.Item(4, 2).Value = "'004404"
$q3Sheet.Cells.Item(4, 3).Value = "平安股息精选沪港深股票C"
$q3Sheet.Cells.Item(4, 4).Value = "'0.02"
$q3Sheet.Cells.Item(4, 5).Value = "'91.74"
$q3Sheet.Cells.Item(4, 6).Value = "'2.41"
$q3Sheet.Cells.Item(4, 7).Value = "'0.0005"
$q3Sheet.Cells.Item(4, 8).Value = 10

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert the new 2022-Q3 row at the top
#    of the data and shift the rest down, appending the 2020-Q4 row that
#    falls off the bottom.
# ---------------------------------------------------------------------------
$totalRows = @(
    @(0, "2022-Q3", 3, 0.27),
    @(1, "2022-Q2", 1, 0.39),
    @(2, "2022-Q1", 1, 0.39),
    @(3, "2021-Q4", 1, 0.39),
    @(4, "2021-Q3", 1, 0.44),
    @(5, "2021-Q2", 1, 0.45),
    @(6, "2021-Q1", 3, 2.48),
    @(7, "2020-Q4", 1, 0.35)
)

# Row 9 is brand new territory for this sheet (it only had rows 1-8 before);
# clone the formatting from row 8 before writing into it.
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 3. Restore "总计" as the active sheet (matches the unedited book view).
# ---------------------------------------------------------------------------
$totalSheet.Activate()
